$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Insert a new row at position 20, pushing existing row 20 (and below) down to row 21
$ws.Rows(20).Insert()

# Fill in the new row 20 with the BITS Jobs attack data
$ws.Range("A20").Value = "T1197-BITS_Jobs_[5].pcapng"
$ws.Range("B20").Value = 1
$ws.Range("C20").Value = "-"
$ws.Range("D20").Value = "15893, 33224"
$ws.Range("E20").Value = "(Interesa más el campo agent=Microsoft BITS 7/8)"
$ws.Range("F20").Value = "El control de aplicaciones detecta el uso de Microsoft BITS para la creación del servicio"

# Match style/formatting (centered, wrap text, Arial 10) used by the rest of the table.
$range = $ws.Range("A20:F20")
$range.Font.Name = "Arial"
$range.Font.Size = 10
$range.HorizontalAlignment = -4108
$range.VerticalAlignment = -4108
$range.WrapText = $true
$ws.Rows(20).RowHeight = $ws.Rows(21).RowHeight

# Restore the sheet selection state left by the editor
$ws.Range("F21").Select()
